$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 607.8
$ws.Range("I33").Value = 563.8
$ws.Range("K33").Value = 563.8
$ws.Range("M33").Value = -334.8

$ws.Range("H62").Value = 3029.5925
$ws.Range("I62").Value = 2895.7917
$ws.Range("K62").Value = 2895.7917
$ws.Range("M62").Value = -2271.7917

$ws.Range("H65").Value = 3029.5925
$ws.Range("I65").Value = 2895.7917
$ws.Range("K65").Value = 14478.9585
$ws.Range("M65").Value = -11358.9585

$ws.Range("H111").Value = 1961.7142
$ws.Range("J111").Value = 1614.5
$ws.Range("L111").Value = 4843.5
$ws.Range("N111").Value = -10977.5

$ws.Range("H125").Value = 1598.8334
$ws.Range("I125").Value = 1788.6
$ws.Range("J125").Value = 650
$ws.Range("K125").Value = 16097.4
$ws.Range("L125").Value = 5850
$ws.Range("M125").Value = -13637.4
$ws.Range("N125").Value = -10770

$ws.Range("H137").Value = 5162.1177
$ws.Range("I137").Value = 5748.15
$ws.Range("K137").Value = 17244.45
$ws.Range("M137").Value = -14694.45

$ws.Range("H141").Value = 7197.2354
$ws.Range("I141").Value = 5181.077
$ws.Range("K141").Value = 15543.231
$ws.Range("M141").Value = -10363.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 50377
$ws.Range("J43").Value = 50377
$ws.Range("L43").Value = 50377
$ws.Range("N43").Value = -51003

$ws.Range("H61").Value = 5711.5386
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 5770.8335
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 5770.8335
$ws.Range("M61").Value = -4788
$ws.Range("N61").Value = -6194.8335

$ws.Range("H74").Value = 1424.579
$ws.Range("I74").Value = 1504.7858
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 1504.7858
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = -630.7858000000001
$ws.Range("N74").Value = -2948

$ws.Range("H77").Value = 1424.579
$ws.Range("I77").Value = 1504.7858
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 7523.929
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = -3155.929
$ws.Range("N77").Value = -14736

$ws.Range("H132").Value = 27029122
$ws.Range("I132").Value = 32259944
$ws.Range("K132").Value = 96779832
$ws.Range("M132").Value = -96777302

$ws.Range("H136").Value = 5711.5386
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 5770.8335
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 17312.5005
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -22412.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3829.8333
$ws.Range("I20").Value = 4996.6665
$ws.Range("J20").Value = 2663
$ws.Range("K20").Value = 4996.6665
$ws.Range("L20").Value = 2663
$ws.Range("M20").Value = -4749.6665
$ws.Range("N20").Value = -3157

$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19516

$ws.Range("H86").Value = 20001296
$ws.Range("I86").Value = 23810688
$ws.Range("J86").Value = 1986.75
$ws.Range("K86").Value = 23810688
$ws.Range("L86").Value = 1986.75
$ws.Range("M86").Value = -23809565
$ws.Range("N86").Value = -4232.75

$ws.Range("H89").Value = 20001296
$ws.Range("I89").Value = 23810688
$ws.Range("J89").Value = 1986.75
$ws.Range("K89").Value = 119053440
$ws.Range("L89").Value = 9933.75
$ws.Range("M89").Value = -119047824
$ws.Range("N89").Value = -21165.75

$ws.Range("H105").Value = 2830.0588
$ws.Range("I105").Value = 1246.909
$ws.Range("K105").Value = 1246.909
$ws.Range("M105").Value = 500.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2290.1875
$ws.Range("I31").Value = 1967.5454
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1967.5454
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1672.5454
$ws.Range("N31").Value = -3590

$ws.Range("H34").Value = 2290.1875
$ws.Range("I34").Value = 1967.5454
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1967.5454
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1765.5454
$ws.Range("N34").Value = -3404

$ws.Range("H41").Value = 18022.357
$ws.Range("J41").Value = 45123
$ws.Range("L41").Value = 45123
$ws.Range("N41").Value = -45979

$ws.Range("H122").Value = 29715440
$ws.Range("J122").Value = 6345.273
$ws.Range("L122").Value = 19035.819
$ws.Range("N122").Value = -23935.819

$ws.Range("H132").Value = 3579.25
$ws.Range("I132").Value = 2356.0833
$ws.Range("K132").Value = 7068.249899999999
$ws.Range("M132").Value = -4538.249899999999

$ws.Range("H134").Value = 2991.8948
$ws.Range("I134").Value = 2288.8
$ws.Range("K134").Value = 6866.400000000001
$ws.Range("M134").Value = -4331.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1158.25
$ws.Range("J25").Value = 1399.8
$ws.Range("L25").Value = 4199.4
$ws.Range("N25").Value = -4537.4

$ws.Range("H30").Value = 1158.25
$ws.Range("J30").Value = 1399.8
$ws.Range("L30").Value = 4199.4
$ws.Range("N30").Value = -4403.4

$ws.Range("H68").Value = 15153870
$ws.Range("I68").Value = 1099
$ws.Range("K68").Value = 3297
$ws.Range("M68").Value = -2486

$ws.Range("H70").Value = 130592.4
$ws.Range("J70").Value = 216820.67
$ws.Range("L70").Value = 650462.01
$ws.Range("N70").Value = -651092.01

$ws.Range("H71").Value = 15153870
$ws.Range("I71").Value = 1099
$ws.Range("K71").Value = 9891
$ws.Range("M71").Value = -5835

$ws.Range("H73").Value = 130592.4
$ws.Range("J73").Value = 216820.67
$ws.Range("L73").Value = 650462.01
$ws.Range("N73").Value = -652646.01

$ws.Range("H80").Value = 7640.3
$ws.Range("J80").Value = 7900.5
$ws.Range("L80").Value = 23701.5
$ws.Range("N80").Value = -25573.5

$ws.Range("H83").Value = 7640.3
$ws.Range("J83").Value = 7900.5
$ws.Range("L83").Value = 71104.5
$ws.Range("N83").Value = -80464.5

$ws.Range("H131").Value = 2293.4285
$ws.Range("J131").Value = 2929
$ws.Range("L131").Value = 8787
$ws.Range("N131").Value = -18867

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 86429.71000000001
$ws.Range("I70").Value = 130393.22
$ws.Range("K70").Value = 130393.22
$ws.Range("M70").Value = -130123.22

$ws.Range("H73").Value = 86429.71000000001
$ws.Range("I73").Value = 130393.22
$ws.Range("K73").Value = 130393.22
$ws.Range("M73").Value = -129457.22

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 700
$ws.Range("I7").Value = 700
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -588
$ws.Range("N7").ClearContents()

$ws.Range("H16").Value = 16666832
$ws.Range("I16").Value = 22727422
$ws.Range("J16").Value = 209.5
$ws.Range("K16").Value = 22727422
$ws.Range("L16").Value = 209.5
$ws.Range("M16").Value = -22727252
$ws.Range("N16").Value = -549.5

$ws.Range("H93").Value = 1563.375
$ws.Range("I93").Value = 1420.6
$ws.Range("J93").Value = 1801.3334
$ws.Range("K93").Value = 1420.6
$ws.Range("L93").Value = 1801.3334
$ws.Range("M93").Value = -172.5999999999999
$ws.Range("N93").Value = -4297.3334

$ws.Range("H122").Value = 58827596
$ws.Range("I122").Value = 142859260
$ws.Range("J122").Value = 5428.2
$ws.Range("K122").Value = 428577780
$ws.Range("L122").Value = 16284.6
$ws.Range("M122").Value = -428575330
$ws.Range("N122").Value = -21184.6

$ws.Range("H126").Value = 700
$ws.Range("I126").Value = 700
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2100
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 370
$ws.Range("N126").ClearContents()

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H107").Value = 472.875
$ws.Range("I107").Value = 197.57143
$ws.Range("K107").Value = 592.71429
$ws.Range("M107").Value = 1327.28571

$ws.Range("H122").Value = 1957.4348
$ws.Range("I122").Value = 1680.7646
$ws.Range("K122").Value = 5042.293799999999
$ws.Range("M122").Value = -2592.293799999999

$ws.Range("H132").Value = 3056.8667
$ws.Range("I132").Value = 2934.6
$ws.Range("K132").Value = 8803.799999999999
$ws.Range("M132").Value = -6273.799999999999

$ws.Range("H138").Value = 59995
$ws.Range("J138").Value = 59995
$ws.Range("L138").Value = 59995
$ws.Range("N138").Value = -70275
